# Lecture 6 and calendar updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calendar2021")

$ws.Range("F5").Value = "HW 3(https://canvas.jmu.edu/courses/1775272/quizzes)"

$ws.Range("F7").Value = "HW4(https://canvas.jmu.edu/courses/1775272/quizzes)"
$ws.Range("G7").Value = "Quiz 0 Retake;"

$ws.Range("D9").Value = "CSPs Part 2; slides(slides/06_CSP_Part2.pdf);video(https://canvas.jmu.edu/courses/1775272/modules)"
$ws.Range("F9").Value = "PA 2"
$ws.Range("G9").Value = "HW4;Quiz 1(mquizzes/mquiz1/mquiz1.php)"

$ws.Range("D10").Value = "Adversarial Search -- Alpha/Beta Pruning"
$ws.Range("E10").Value = "Chp 5.1 - 5.4"

$ws.Range("D11").Value = "Heuristic Alpha/Beta Search, Monte Carlo Search, Chance Trees"
$ws.Range("E11").Value = "Chp 5.5 - 5.7"

$ws.Range("D12").Value = "Chance Tree Lab"

$ws.Activate()
$ws.Range("E12").Select()
